$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bioSampleNumber values in column C (rows 2-19): add 100 to each existing value
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $newValue = ($row - 1) + 100
    $cell.Value = $newValue
}

# Adjust row height for data rows (2-19) from 15 to 13.8
$dataRows = $ws.Range("A2:A19").EntireRow
$dataRows.RowHeight = 13.8

# Move the active selection to G20
$ws.Range("G20").Select()
